$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 3409.1428
$ws.Range("I32").Value = 1130
$ws.Range("K32").Value = 1130
$ws.Range("M32").Value = -804

# Row 47
$ws.Range("H47").Value = 9235.666999999999
$ws.Range("I47").Value = 1353.5
$ws.Range("J47").Value = 25000
$ws.Range("K47").Value = 1353.5
$ws.Range("L47").Value = 25000
$ws.Range("M47").Value = -381.5
$ws.Range("N47").Value = -26944

# Row 58
$ws.Range("H58").Value = 4344.7646
$ws.Range("J58").Value = 14393
$ws.Range("L58").Value = 43179
$ws.Range("N58").Value = -43479

# Row 59
$ws.Range("H59").Value = 10000
$ws.Range("J59").Value = 10000
$ws.Range("L59").Value = 30000
$ws.Range("N59").Value = -31114

# Row 76
$ws.Range("H76").Value = 21733.166
$ws.Range("I76").Value = 21166.334
$ws.Range("J76").Value = 22300
$ws.Range("K76").Value = 21166.334
$ws.Range("L76").Value = 22300
$ws.Range("M76").Value = -20851.334
$ws.Range("N76").Value = -22930

# Row 79
$ws.Range("H79").Value = 21733.166
$ws.Range("I79").Value = 21166.334
$ws.Range("J79").Value = 22300
$ws.Range("K79").Value = 21166.334
$ws.Range("L79").Value = 22300
$ws.Range("M79").Value = -20074.334
$ws.Range("N79").Value = -24484

# Row 103
$ws.Range("H103").Value = 395.45
$ws.Range("I103").Value = 522.4286
$ws.Range("J103").Value = 327.07693
$ws.Range("K103").Value = 1567.2858
$ws.Range("L103").Value = 981.2307900000001
$ws.Range("M103").Value = -981.2857999999999
$ws.Range("N103").Value = -2153.23079

$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 3295.7693
$ws.Range("I88").Value = 2966.3333
$ws.Range("J88").Value = 3394.6
$ws.Range("K88").Value = 2966.3333
$ws.Range("L88").Value = 3394.6
$ws.Range("M88").Value = -2560.3333
$ws.Range("N88").Value = -4206.6

# Row 91
$ws.Range("H91").Value = 3295.7693
$ws.Range("I91").Value = 2966.3333
$ws.Range("J91").Value = 3394.6
$ws.Range("K91").Value = 2966.3333
$ws.Range("L91").Value = 3394.6
$ws.Range("M91").Value = -1562.3333
$ws.Range("N91").Value = -6202.6

# Row 97
$ws.Range("H97").Value = 475.33334
$ws.Range("I97").Value = 249.125
$ws.Range("J97").Value = 1199.2
$ws.Range("K97").Value = 249.125
$ws.Range("L97").Value = 1199.2
$ws.Range("M97").Value = 246.875
$ws.Range("N97").Value = -2191.2

# Row 122
$ws.Range("H122").Value = 2799.7646
$ws.Range("I122").Value = 2495.25
$ws.Range("J122").Value = 4713.857
$ws.Range("K122").Value = 7485.75
$ws.Range("L122").Value = 14141.571
$ws.Range("M122").Value = -5035.75
$ws.Range("N122").Value = -19041.571

$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 34647.332
$ws.Range("I26").Value = 34647.332
$ws.Range("K26").Value = 34647.332
$ws.Range("M26").Value = -34355.332

# Row 107
$ws.Range("H107").Value = 4408.7646
$ws.Range("J107").Value = 5794.5386
$ws.Range("L107").Value = 5794.5386
$ws.Range("N107").Value = -9634.5386

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1422.3077
$ws.Range("I31").Value = 1502
$ws.Range("J31").Value = 984
$ws.Range("K31").Value = 1502
$ws.Range("L31").Value = 984
$ws.Range("M31").Value = -1207
$ws.Range("N31").Value = -1574

# Row 34
$ws.Range("H34").Value = 1422.3077
$ws.Range("I34").Value = 1502
$ws.Range("J34").Value = 984
$ws.Range("K34").Value = 1502
$ws.Range("L34").Value = 984
$ws.Range("M34").Value = -1300
$ws.Range("N34").Value = -1388

# Row 107
$ws.Range("H107").Value = 958.34283
$ws.Range("I107").Value = 735
$ws.Range("J107").Value = 987.1613
$ws.Range("K107").Value = 735
$ws.Range("L107").Value = 987.1613
$ws.Range("M107").Value = 1185
$ws.Range("N107").Value = -4827.1613

$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Range("H32").Value = 335199.66
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

# Row 35
$ws.Range("H35").Value = 200
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# Row 38
$ws.Range("H38").Value = 139.3077
$ws.Range("I38").Value = 63.8
$ws.Range("J38").Value = 186.5
$ws.Range("K38").Value = 191.4
$ws.Range("L38").Value = 559.5
$ws.Range("M38").Value = 155.6
$ws.Range("N38").Value = -1253.5

# Row 41
$ws.Range("H41").Value = 62
$ws.Range("I41").Value = 99
$ws.Range("K41").Value = 297
$ws.Range("M41").Value = 41

# Row 42
$ws.Range("H42").Value = 1061
$ws.Range("J42").Value = 1061
$ws.Range("L42").Value = 3183
$ws.Range("N42").Value = -4251

# Row 49
$ws.Range("H49").Value = 4662.5
$ws.Range("I49").Value = 2825
$ws.Range("K49").Value = 8475
$ws.Range("M49").Value = -8319

# Row 57
$ws.Range("H57").Value = 211701.8
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

# Row 59
$ws.Range("H59").Value = 6750
$ws.Range("I59").Value = 6750
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 20250
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -19710

# Row 69
$ws.Range("H69").Value = 6111.1113
$ws.Range("I69").Value = 6687.5
$ws.Range("J69").Value = 1500
$ws.Range("K69").Value = 20062.5
$ws.Range("L69").Value = 4500
$ws.Range("M69").Value = -19251.5
$ws.Range("N69").Value = -6122

# Row 72
$ws.Range("H72").Value = 6111.1113
$ws.Range("I72").Value = 6687.5
$ws.Range("J72").Value = 1500
$ws.Range("K72").Value = 60187.5
$ws.Range("L72").Value = 13500
$ws.Range("M72").Value = -56131.5
$ws.Range("N72").Value = -21612

# Row 80
$ws.Range("H80").Value = 5997.125
$ws.Range("J80").Value = 4749.5
$ws.Range("L80").Value = 14248.5
$ws.Range("N80").Value = -16120.5

# Row 82
$ws.Range("H82").Value = 12496.5
$ws.Range("I82").Value = 4993
$ws.Range("K82").Value = 14979
$ws.Range("M82").Value = -14573

# Row 83
$ws.Range("H83").Value = 5997.125
$ws.Range("J83").Value = 4749.5
$ws.Range("L83").Value = 42745.5
$ws.Range("N83").Value = -52105.5

# Row 85
$ws.Range("H85").Value = 12496.5
$ws.Range("I85").Value = 4993
$ws.Range("K85").Value = 14979
$ws.Range("M85").Value = -13575

# Row 87
$ws.Range("H87").Value = 40000
$ws.Range("I87").Value = 40000
$ws.Range("K87").Value = 120000
$ws.Range("M87").Value = -118752

# Row 90
$ws.Range("H90").Value = 40000
$ws.Range("I90").Value = 40000
$ws.Range("K90").Value = 360000
$ws.Range("M90").Value = -353760

$ws = $wb.Worksheets.Item("GSM")
# Row 53
$ws.Range("H53").Value = 20000
$ws.Range("I53").Value = 10000
$ws.Range("J53").Value = 30000
$ws.Range("K53").Value = 10000
$ws.Range("L53").Value = 30000
$ws.Range("M53").Value = -9369
$ws.Range("N53").Value = -31262

# Row 102
$ws.Range("H102").Value = 2322.889
$ws.Range("I102").Value = 2110.3333
$ws.Range("J102").Value = 2748
$ws.Range("K102").Value = 2110.3333
$ws.Range("L102").Value = 2748
$ws.Range("M102").Value = -488.3332999999998
$ws.Range("N102").Value = -5992

# Row 126
$ws.Range("H126").Value = 2843.85
$ws.Range("I126").Value = 3341.9092
$ws.Range("J126").Value = 2235.111
$ws.Range("K126").Value = 10025.7276
$ws.Range("L126").Value = 6705.333
$ws.Range("M126").Value = -7555.7276
$ws.Range("N126").Value = -11645.333

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1342.2858
$ws.Range("J46").Value = 1342.2858
$ws.Range("L46").Value = 1342.2858
$ws.Range("N46").Value = -1718.2858

# Row 55
$ws.Range("H55").Value = 890.06665
$ws.Range("I55").Value = 481.77777
$ws.Range("J55").Value = 1502.5
$ws.Range("K55").Value = 481.77777
$ws.Range("L55").Value = 1502.5
$ws.Range("M55").Value = -308.77777
$ws.Range("N55").Value = -1848.5

$ws = $wb.Worksheets.Item("WVR")
# Row 76
$ws.Range("H76").Value = 23499.5
$ws.Range("J76").Value = 23499.5
$ws.Range("L76").Value = 23499.5
$ws.Range("N76").Value = -24129.5

# Row 79
$ws.Range("H79").Value = 23499.5
$ws.Range("J79").Value = 23499.5
$ws.Range("L79").Value = 23499.5
$ws.Range("N79").Value = -25683.5
